# Ajusta controle de investimento
# The "Quanto investir por mês?" value is updated from 2440 to the
# suggested investment amount (Salario*30% = 3030). All the other
# cells on the "APP" sheet (patrimônio acumulado, dividendos, cenários
# and the FII allocation table) are driven by formulas that reference
# this cell (directly or via the named range "Valor_investir_mês"), so
# they recalculate automatically once the value below is changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("APP")

$ws.Range("D17").Value = 3030

$excel.Calculate()
